$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-01-14 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-01-15 Sunday", 2) | Out-Null
$d.Content.Find.Execute("61+15=", $true, $false, $false, $false, $false, $true, 1, $false, "81+0=", 2) | Out-Null
$d.Content.Find.Execute("63-41=", $true, $false, $false, $false, $false, $true, 1, $false, "50+38=", 2) | Out-Null
$d.Content.Find.Execute("67+3=", $true, $false, $false, $false, $false, $true, 1, $false, "52+14=", 2) | Out-Null
$d.Content.Find.Execute("10-0=", $true, $false, $false, $false, $false, $true, 1, $false, "18-18=", 2) | Out-Null
$d.Content.Find.Execute("39-34=", $true, $false, $false, $false, $false, $true, 1, $false, "52-20=", 2) | Out-Null
$d.Content.Find.Execute("20+79=", $true, $false, $false, $false, $false, $true, 1, $false, "58+11=", 2) | Out-Null
$d.Content.Find.Execute("26+16=", $true, $false, $false, $false, $false, $true, 1, $false, "78-19=", 2) | Out-Null
$d.Content.Find.Execute("71-27=", $true, $false, $false, $false, $false, $true, 1, $false, "56-48=", 2) | Out-Null
$d.Content.Find.Execute("36-24=", $true, $false, $false, $false, $false, $true, 1, $false, "63-24=", 2) | Out-Null
$d.Content.Find.Execute("48+51=", $true, $false, $false, $false, $false, $true, 1, $false, "51+28=", 2) | Out-Null
$d.Content.Find.Execute("11+50=", $true, $false, $false, $false, $false, $true, 1, $false, "99-57=", 2) | Out-Null
$d.Content.Find.Execute("23-3=", $true, $false, $false, $false, $false, $true, 1, $false, "77-28=", 2) | Out-Null
$d.Content.Find.Execute("56+31=", $true, $false, $false, $false, $false, $true, 1, $false, "45+30=", 2) | Out-Null
$d.Content.Find.Execute("17+65=", $true, $false, $false, $false, $false, $true, 1, $false, "13+34=", 2) | Out-Null
$d.Content.Find.Execute("88-58=", $true, $false, $false, $false, $false, $true, 1, $false, "10+24=", 2) | Out-Null
$d.Content.Find.Execute("91+3=", $true, $false, $false, $false, $false, $true, 1, $false, "41-5=", 2) | Out-Null
$d.Content.Find.Execute("77-59=", $true, $false, $false, $false, $false, $true, 1, $false, "58-10=", 2) | Out-Null
$d.Content.Find.Execute("85-36=", $true, $false, $false, $false, $false, $true, 1, $false, "74-47=", 2) | Out-Null
$d.Content.Find.Execute("40+32=", $true, $false, $false, $false, $false, $true, 1, $false, "14+59=", 2) | Out-Null
$d.Content.Find.Execute("31+21=", $true, $false, $false, $false, $false, $true, 1, $false, "40+34=", 2) | Out-Null
$d.Content.Find.Execute("77-54=", $true, $false, $false, $false, $false, $true, 1, $false, "98-34=", 2) | Out-Null
$d.Content.Find.Execute("68-34=", $true, $false, $false, $false, $false, $true, 1, $false, "67-32=", 2) | Out-Null
$d.Content.Find.Execute("75-44=", $true, $false, $false, $false, $false, $true, 1, $false, "83-52=", 2) | Out-Null
$d.Content.Find.Execute("50+39=", $true, $false, $false, $false, $false, $true, 1, $false, "45-5=", 2) | Out-Null
$d.Content.Find.Execute("46+52=", $true, $false, $false, $false, $false, $true, 1, $false, "88-2=", 2) | Out-Null
$d.Content.Find.Execute("27+12=", $true, $false, $false, $false, $false, $true, 1, $false, "58-28=", 2) | Out-Null
$d.Content.Find.Execute("19+10=", $true, $false, $false, $false, $false, $true, 1, $false, "59+30=", 2) | Out-Null
$d.Content.Find.Execute("96-15=", $true, $false, $false, $false, $false, $true, 1, $false, "85-60=", 2) | Out-Null
$d.Content.Find.Execute("18-5=", $true, $false, $false, $false, $false, $true, 1, $false, "67+4=", 2) | Out-Null
$d.Content.Find.Execute("38-12=", $true, $false, $false, $false, $false, $true, 1, $false, "73-29=", 2) | Out-Null
$d.Content.Find.Execute("48+15=", $true, $false, $false, $false, $false, $true, 1, $false, "1+15=", 2) | Out-Null
$d.Content.Find.Execute("33+18=", $true, $false, $false, $false, $false, $true, 1, $false, "51+18=", 2) | Out-Null
$d.Content.Find.Execute("61+8=", $true, $false, $false, $false, $false, $true, 1, $false, "99-46=", 2) | Out-Null
$d.Content.Find.Execute("25+35=", $true, $false, $false, $false, $false, $true, 1, $false, "73-33=", 2) | Out-Null
$d.Content.Find.Execute("14+12=", $true, $false, $false, $false, $false, $true, 1, $false, "50-22=", 2) | Out-Null
$d.Content.Find.Execute("84+8=", $true, $false, $false, $false, $false, $true, 1, $false, "4+86=", 2) | Out-Null
$d.Content.Find.Execute("30-6=", $true, $false, $false, $false, $false, $true, 1, $false, "72-28=", 2) | Out-Null
$d.Content.Find.Execute("75+21=", $true, $false, $false, $false, $false, $true, 1, $false, "90-56=", 2) | Out-Null
$d.Content.Find.Execute("17-14=", $true, $false, $false, $false, $false, $true, 1, $false, "81-1=", 2) | Out-Null
$d.Content.Find.Execute("73-56=", $true, $false, $false, $false, $false, $true, 1, $false, "51+22=", 2) | Out-Null
$d.Content.Find.Execute("29+20=", $true, $false, $false, $false, $false, $true, 1, $false, "1+31=", 2) | Out-Null
$d.Content.Find.Execute("14+37=", $true, $false, $false, $false, $false, $true, 1, $false, "28+47=", 2) | Out-Null
$d.Content.Find.Execute("75-6=", $true, $false, $false, $false, $false, $true, 1, $false, "44+39=", 2) | Out-Null
$d.Content.Find.Execute("92-7=", $true, $false, $false, $false, $false, $true, 1, $false, "27+42=", 2) | Out-Null
$d.Content.Find.Execute("84-80=", $true, $false, $false, $false, $false, $true, 1, $false, "26+65=", 2) | Out-Null
$d.Content.Find.Execute("66-34=", $true, $false, $false, $false, $false, $true, 1, $false, "57+3=", 2) | Out-Null
$d.Content.Find.Execute("42-38=", $true, $false, $false, $false, $false, $true, 1, $false, "83-57=", 2) | Out-Null
$d.Content.Find.Execute("45-6=", $true, $false, $false, $false, $false, $true, 1, $false, "8-2=", 2) | Out-Null
$d.Content.Find.Execute("32-21=", $true, $false, $false, $false, $false, $true, 1, $false, "0+69=", 2) | Out-Null
$d.Content.Find.Execute("21+54=", $true, $false, $false, $false, $false, $true, 1, $false, "20+68=", 2) | Out-Null
$d.Content.Find.Execute("76-63=", $true, $false, $false, $false, $false, $true, 1, $false, "65-22=", 2) | Out-Null
$d.Content.Find.Execute("26+6=", $true, $false, $false, $false, $false, $true, 1, $false, "27+4=", 2) | Out-Null
$d.Content.Find.Execute("95-68=", $true, $false, $false, $false, $false, $true, 1, $false, "5+41=", 2) | Out-Null
$d.Content.Find.Execute("24+73=", $true, $false, $false, $false, $false, $true, 1, $false, "50-40=", 2) | Out-Null
$d.Content.Find.Execute("42+14=", $true, $false, $false, $false, $false, $true, 1, $false, "43+37=", 2) | Out-Null
$d.Content.Find.Execute("59+26=", $true, $false, $false, $false, $false, $true, 1, $false, "18+3=", 2) | Out-Null
$d.Content.Find.Execute("28+12=", $true, $false, $false, $false, $false, $true, 1, $false, "36-32=", 2) | Out-Null
$d.Content.Find.Execute("72-69=", $true, $false, $false, $false, $false, $true, 1, $false, "59-26=", 2) | Out-Null
$d.Content.Find.Execute("74+15=", $true, $false, $false, $false, $false, $true, 1, $false, "62-10=", 2) | Out-Null
$d.Content.Find.Execute("66+16=", $true, $false, $false, $false, $false, $true, 1, $false, "98-13=", 2) | Out-Null
$d.Content.Find.Execute("98-37=", $true, $false, $false, $false, $false, $true, 1, $false, "84-37=", 2) | Out-Null
$d.Content.Find.Execute("31-12=", $true, $false, $false, $false, $false, $true, 1, $false, "95-24=", 2) | Out-Null
$d.Content.Find.Execute("10+66=", $true, $false, $false, $false, $false, $true, 1, $false, "32+12=", 2) | Out-Null
$d.Content.Find.Execute("4+13=", $true, $false, $false, $false, $false, $true, 1, $false, "78-75=", 2) | Out-Null
$d.Content.Find.Execute("72+19=", $true, $false, $false, $false, $false, $true, 1, $false, "32+32=", 2) | Out-Null
$d.Content.Find.Execute("19-3=", $true, $false, $false, $false, $false, $true, 1, $false, "15+19=", 2) | Out-Null
$d.Content.Find.Execute("66+24=", $true, $false, $false, $false, $false, $true, 1, $false, "68-10=", 2) | Out-Null
$d.Content.Find.Execute("66+29=", $true, $false, $false, $false, $false, $true, 1, $false, "24-18=", 2) | Out-Null
$d.Content.Find.Execute("95-59=", $true, $false, $false, $false, $false, $true, 1, $false, "52+34=", 2) | Out-Null
$d.Content.Find.Execute("21-12=", $true, $false, $false, $false, $false, $true, 1, $false, "33-32=", 2) | Out-Null
$d.Content.Find.Execute("3+4=", $true, $false, $false, $false, $false, $true, 1, $false, "31-25=", 2) | Out-Null
$d.Content.Find.Execute("58-34=", $true, $false, $false, $false, $false, $true, 1, $false, "88-70=", 2) | Out-Null
$d.Content.Find.Execute("12+56=", $true, $false, $false, $false, $false, $true, 1, $false, "3+81=", 2) | Out-Null
$d.Content.Find.Execute("75+5=", $true, $false, $false, $false, $false, $true, 1, $false, "53+41=", 2) | Out-Null
$d.Content.Find.Execute("37+0=", $true, $false, $false, $false, $false, $true, 1, $false, "48-46=", 2) | Out-Null
$d.Content.Find.Execute("92-51=", $true, $false, $false, $false, $false, $true, 1, $false, "55+16=", 2) | Out-Null
$d.Content.Find.Execute("74-14=", $true, $false, $false, $false, $false, $true, 1, $false, "32-11=", 2) | Out-Null
$d.Content.Find.Execute("63+15=", $true, $false, $false, $false, $false, $true, 1, $false, "87+2=", 2) | Out-Null
$d.Content.Find.Execute("95-0=", $true, $false, $false, $false, $false, $true, 1, $false, "19-18=", 2) | Out-Null
$d.Content.Find.Execute("62-27=", $true, $false, $false, $false, $false, $true, 1, $false, "42-7=", 2) | Out-Null
$d.Content.Find.Execute("72+12=", $true, $false, $false, $false, $false, $true, 1, $false, "50+35=", 2) | Out-Null
$d.Content.Find.Execute("77-13=", $true, $false, $false, $false, $false, $true, 1, $false, "14+54=", 2) | Out-Null
$d.Content.Find.Execute("9-7=", $true, $false, $false, $false, $false, $true, 1, $false, "93-72=", 2) | Out-Null
$d.Content.Find.Execute("49+21=", $true, $false, $false, $false, $false, $true, 1, $false, "12-11=", 2) | Out-Null
$d.Content.Find.Execute("64+21=", $true, $false, $false, $false, $false, $true, 1, $false, "74-50=", 2) | Out-Null
$d.Content.Find.Execute("7+9=", $true, $false, $false, $false, $false, $true, 1, $false, "75-61=", 2) | Out-Null
$d.Content.Find.Execute("47+37=", $true, $false, $false, $false, $false, $true, 1, $false, "17+52=", 2) | Out-Null
$d.Content.Find.Execute("85-38=", $true, $false, $false, $false, $false, $true, 1, $false, "4+79=", 2) | Out-Null
$d.Content.Find.Execute("30+59=", $true, $false, $false, $false, $false, $true, 1, $false, "34-26=", 2) | Out-Null
$d.Content.Find.Execute("81+11=", $true, $false, $false, $false, $false, $true, 1, $false, "82+6=", 2) | Out-Null
$d.Content.Find.Execute("82-39=", $true, $false, $false, $false, $false, $true, 1, $false, "19+19=", 2) | Out-Null
$d.Content.Find.Execute("63-48=", $true, $false, $false, $false, $false, $true, 1, $false, "67-18=", 2) | Out-Null
$d.Content.Find.Execute("7+13=", $true, $false, $false, $false, $false, $true, 1, $false, "43+25=", 2) | Out-Null
$d.Content.Find.Execute("19-8=", $true, $false, $false, $false, $false, $true, 1, $false, "50-6=", 2) | Out-Null
$d.Content.Find.Execute("16+9=", $true, $false, $false, $false, $false, $true, 1, $false, "39+7=", 2) | Out-Null
$d.Content.Find.Execute("96-43=", $true, $false, $false, $false, $false, $true, 1, $false, "72-49=", 2) | Out-Null
$d.Content.Find.Execute("44+25=", $true, $false, $false, $false, $false, $true, 1, $false, "80+14=", 2) | Out-Null
$d.Content.Find.Execute("36-15=", $true, $false, $false, $false, $false, $true, 1, $false, "15+3=", 2) | Out-Null
$d.Content.Find.Execute("6+20=", $true, $false, $false, $false, $false, $true, 1, $false, "33-29=", 2) | Out-Null
$d.Content.Find.Execute("85-50=", $true, $false, $false, $false, $false, $true, 1, $false, "98-55=", 2) | Out-Null
